$d = $word.ActiveDocument
$n = $d.Paragraphs.Count

# Locate the anchor paragraph ("LOB1053: ...") and the trailing footer block
# ("Ver no Jupiter ..." / "(c) 2020 ...") that follows it, then remove the
# blank line + those two footer paragraphs, leaving the single blank
# paragraph that originally sat right before the page-break paragraph.
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*LOB1053: Física III (Requisito fraco)*") {
        $startIdx = $i + 1
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endIdx = $i
    }
}

if ($startIdx -ne -1 -and $endIdx -ne -1 -and $endIdx -ge $startIdx) {
    $start = $d.Paragraphs($startIdx).Range.Start
    $end = $d.Paragraphs($endIdx).Range.End
    $delRange = $d.Range($start, $end)
    $delRange.Delete()
    Write-Output "Removed paragraphs $startIdx..$endIdx"
} else {
    Write-Output "Target paragraphs not found; no change made."
}
